# NewlineInFormulas.xlsx: add a Fibonacci-style shared-formula column on
# Sheet1 (B1:B10) next to the existing multi-line SUM() formula in A1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Seed values
$ws.Range("B1").Value = 1
$ws.Range("B2").Value = 2

# B3:B10 = previous two cells added together, entered once across the
# whole range so Excel stores it as a single shared formula (t="shared").
$ws.Range("B3:B10").Formula = "=B1+B2"

# Matches the selection recorded in the saved sheetView (B3 active cell,
# B3:B10 selected) after filling the formula down.
$ws.Range("B3:B10").Select() | Out-Null
